$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite columns A and B with the consolidated Capability/Enabled data
$ws.Cells.Item(1,1).Value = "Capability Name"
$ws.Cells.Item(1,2).Value = "Enabled (YES/NO)"
$ws.Cells.Item(2,1).Value = "Case Management"
$ws.Cells.Item(2,2).Value = "YES"
$ws.Cells.Item(3,1).Value = "Case Assignment Rules"
$ws.Cells.Item(3,2).Value = "YES"
$ws.Cells.Item(4,1).Value = "Escalation Rules"
$ws.Cells.Item(4,2).Value = "NO"
$ws.Cells.Item(5,1).Value = "Auto Response Rules"
$ws.Cells.Item(5,2).Value = "NO"
$ws.Cells.Item(6,1).Value = "Omni"
$ws.Cells.Item(6,2).Value = "NO"
$ws.Cells.Item(7,1).Value = "Entitlements"
$ws.Cells.Item(7,2).Value = "YES"
$ws.Cells.Item(8,1).Value = "Milestones"
$ws.Cells.Item(8,2).Value = "NO"
$ws.Cells.Item(9,1).Value = "Knowledge Base"
$ws.Cells.Item(9,2).Value = "NO"
$ws.Cells.Item(10,1).Value = "Email"
$ws.Cells.Item(10,2).Value = "NO"
$ws.Cells.Item(11,1).Value = "Web"
$ws.Cells.Item(11,2).Value = "NO"
$ws.Cells.Item(12,1).Value = "Macros"
$ws.Cells.Item(12,2).Value = "NO"
$ws.Cells.Item(13,1).Value = "Quick Text"
$ws.Cells.Item(13,2).Value = "NO"
$ws.Cells.Item(14,1).Value = "Service Console"
$ws.Cells.Item(14,2).Value = "NO"
$ws.Cells.Item(15,1).Value = "Live Chat"
$ws.Cells.Item(15,2).Value = "NO"
$ws.Cells.Item(16,1).Value = "CTI Integration"
$ws.Cells.Item(16,2).Value = "NO"
$ws.Cells.Item(17,1).Value = "Field Service Management"
$ws.Cells.Item(17,2).Value = "NO"
$ws.Cells.Item(18,1).Value = "Einstein Bots"
$ws.Cells.Item(18,2).Value = "NO"
$ws.Cells.Item(19,1).Value = "Case Queues"
$ws.Cells.Item(19,2).Value = "YES"
$ws.Cells.Item(20,1).Value = "Service Level Agreements(SLAs)"
$ws.Cells.Item(20,2).Value = "YES"
$ws.Cells.Item(21,1).Value = "Reports"
$ws.Cells.Item(21,2).Value = "YES"
$ws.Cells.Item(22,1).Value = "Dashboards"
$ws.Cells.Item(22,2).Value = "YES"
$ws.Cells.Item(23,1).Value = "Automation Flows"
$ws.Cells.Item(23,2).Value = "NO"
$ws.Cells.Item(24,1).Value = "Apex Triggers"
$ws.Cells.Item(24,2).Value = "NO"
$ws.Cells.Item(25,1).Value = "API Integrations"
$ws.Cells.Item(25,2).Value = "NO"
$ws.Cells.Item(26,1).Value = "Data Storage Monitoring"
$ws.Cells.Item(26,2).Value = "NO"
$ws.Cells.Item(27,1).Value = "Permission Sets"
$ws.Cells.Item(27,2).Value = "NO"
$ws.Cells.Item(28,1).Value = "Profiles"
$ws.Cells.Item(28,2).Value = "NO"
$ws.Cells.Item(29,1).Value = "Service Contracts"
$ws.Cells.Item(29,2).Value = "NO"
$ws.Cells.Item(30,1).Value = "Customer Community"
$ws.Cells.Item(30,2).Value = "NO"
$ws.Cells.Item(31,1).Value = "Partner Community"
$ws.Cells.Item(31,2).Value = "NO"
$ws.Cells.Item(32,1).Value = "Incident Management"
$ws.Cells.Item(32,2).Value = "NO"
$ws.Cells.Item(33,1).Value = "Problem Management"
$ws.Cells.Item(33,2).Value = "NO"
$ws.Cells.Item(34,1).Value = "Change Management"
$ws.Cells.Item(34,2).Value = "NO"
$ws.Cells.Item(35,1).Value = "Asset Management"
$ws.Cells.Item(35,2).Value = "NO"
$ws.Cells.Item(36,1).Value = "Work Orders"
$ws.Cells.Item(36,2).Value = "NO"
$ws.Cells.Item(37,1).Value = "Service Appointments"
$ws.Cells.Item(37,2).Value = "NO"
$ws.Cells.Item(38,1).Value = "Routing Configurations"
$ws.Cells.Item(38,2).Value = "NO"
$ws.Cells.Item(39,1).Value = "Skills"
$ws.Cells.Item(39,2).Value = "NO"
$ws.Cells.Item(40,1).Value = "Presence Configurations"
$ws.Cells.Item(40,2).Value = "NO"
$ws.Cells.Item(41,1).Value = "Social Customer Service"
$ws.Cells.Item(41,2).Value = "NO"
$ws.Cells.Item(42,1).Value = "Messaging Channels"
$ws.Cells.Item(42,2).Value = "NO"
$ws.Cells.Item(43,1).Value = "Case Swarming"
$ws.Cells.Item(43,2).Value = "YES"
$ws.Cells.Item(44,1).Value = "Service Analytics"
$ws.Cells.Item(44,2).Value = "NO"
$ws.Cells.Item(45,1).Value = "Customer Feedback Surveys"
$ws.Cells.Item(45,2).Value = "NO"
$ws.Cells.Item(46,1).Value = "Knowledge Recommendations"
$ws.Cells.Item(46,2).Value = "NO"
$ws.Cells.Item(47,1).Value = "Einstein Case Classification"
$ws.Cells.Item(47,2).Value = "YES"
$ws.Cells.Item(48,1).Value = "Einstein Next Best Action"
$ws.Cells.Item(48,2).Value = "NO"
$ws.Cells.Item(49,1).Value = "Service Forecasting"
$ws.Cells.Item(49,2).Value = "NO"

# Remove the old "Description" column entirely (content now consolidated into A/B)
$ws.Columns.Item(3).Delete()

# Drop the two now-unused trailing rows so the used range shrinks to A1:B49
$ws.Rows.Item(51).Delete()
$ws.Rows.Item(50).Delete()

# Style the header row: bold font, thin box border, centered horizontally, top vertically
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
